# edit.ps1 - applies the "Divided votes from voters" edit to fiszka.docx
# via Word COM-interop (headless iron_native runtime).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    # NOTE: MatchCase must stay $false - this runtime's Find mis-handles
    # case-sensitive matching together with accented (non-ASCII) text.
    $ok = $range.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: could not find text: $old"
    }
    return $ok
}

# 1. "Przyspieszenie procesu liczenia głosów " -> drop trailing space
Replace-Text "Przyspieszenie procesu liczenia głosów " "Przyspieszenie procesu liczenia głosów"

# 2. Append a trailing space run after the "...obwodowej komisji wyborczej" bullet
Replace-Text "możliwienie wyborcom wzięcia udziału w głosowaniu bez konieczności dojazdu do obwodowej komisji wyborczej" "możliwienie wyborcom wzięcia udziału w głosowaniu bez konieczności dojazdu do obwodowej komisji wyborczej "

# 3. Week 2 - drop "oraz ustalenie wykorzystywanych narzędzi do ich wdrożenia"
Replace-Text "Tydzień 2. – Projekt graficzny interfejsów graficznych aplikacji oraz ustalenie wykorzystywanych narzędzi do ich wdrożenia. Zaimplementowanie bazy danych. " "Tydzień 2. – Projekt graficzny interfejsów graficznych aplikacji. Zaimplementowanie bazy danych. "

# 4. Week 3 - drop "formularzy oraz "
Replace-Text "Implementacja formularzy oraz podstawowych funkcji logowania użytkowników i ich rejestracji." "Implementacja podstawowych funkcji logowania użytkowników i ich rejestracji."

# 5. Week 4 - "Kontynuacja implementacji oraz początek wdrażania funkcji związanych z wyborami." becomes "Implementacja formularzy."
Replace-Text "Tydzień 4. – Kontynuacja implementacji oraz początek wdrażania funkcji związanych z wyborami. " "Tydzień 4. – Implementacja formularzy. "

# 6. Week 5 - "Dalsza część implementacji." becomes "Implementacja podstawowych funkcji wyborczych."
Replace-Text "Tydzień 5. – Dalsza część implementacji. " "Tydzień 5. – Implementacja podstawowych funkcji wyborczych. "

# 7. Week 6 - "Tydzień 6. - Ostatni etap implementacji." becomes "Tydzień 6. – Wdrożenie zaawansowanych czynności dotyczących przeprowadzania głosowania."
Replace-Text "Tydzień 6. - Ostatni etap implementacji. " "Tydzień 6. – Wdrożenie zaawansowanych czynności dotyczących przeprowadzania głosowania. "

# 8. Week 7 - drop "Wprowadzanie poprawek w kodzie na podstawie wykonanych testów.  " and simplify "Sprawdzenie i ewentualne korygowanie" -> "Sprawdzenie"
Replace-Text "Tydzień 7. – Testowanie oprogramowania. Wprowadzanie poprawek w kodzie na podstawie wykonanych testów.  Sprawdzenie i ewentualne korygowanie zabezpieczeń oprogramowania." "Tydzień 7. – Testowanie oprogramowania. Sprawdzenie zabezpieczeń oprogramowania."

# 9. Week 8 - "Dalsza część testów i prac nad zabezpieczeniami." becomes "Wprowadzenie poprawek w kodzie i zabezpieczeniach na podstawie wykonanych testów.  "
Replace-Text "Tydzień 8. – Dalsza część testów i prac nad zabezpieczeniami." "Tydzień 8. – Wprowadzenie poprawek w kodzie i zabezpieczeniach na podstawie wykonanych testów.  "

# 10. Insert an empty paragraph between "Tydzień 10. ..." and "Kluczowe ryzyka"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Tydzie*10*Projekt gotowy do oddania*") {
        $p.Range.InsertParagraphAfter()
        break
    }
}

# Re-locate the freshly inserted blank paragraph (the one right before
# "Kluczowe ryzyka") and strip its list numbering / reset indentation so it
# renders as a plain indented blank line instead of a bullet.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $paras.Item($i)
    if ($para.Range.Text -like "*Kluczowe ryzyka*") {
        $blank = $paras.Item($i - 1)
        $blank.Range.ListFormat.RemoveNumbers()
        $blank.Range.ParagraphFormat.LeftIndent = 72
        break
    }
}

# 11. "Długotrwała choroba członka zespołu" - text unchanged (run-merge only), no-op

# 12. "autentykacji" -> "uwierzytelniania"
Replace-Text "Profilu Zaufanego do autentykacji użytkowników" "Profilu Zaufanego do uwierzytelniania użytkowników"

Write-Host "edits applied"
